$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DMZ (reverse-)proxy requirement rows -- "dmz reverse proxy; TOGAF docs"
$newRows = @(
    @("External web sites requested from internal web clients will pass through a DMZ web proxy server", "Proxy"),
    @("Blacklisted external web sites will not be accessible from internal web clients", "Proxy"),
    @("The proxy server will scan web traffic for malware.", "Proxy"),
    @("The proxy server will block web traffic when malware is detected.", "Proxy"),
    @("Internal web sites requested from external  web clients will pass through a DMZ reverse web proxy server", "Reverse Proxy"),
    @("The reverse proxy server will scan web traffic for malware.", "Reverse Proxy"),
    @("Blacklisted external ip address will not be able to access the internal web server", "Reverse Proxy"),
    @("The reverse proxy server will block web traffic when malware is detected.", "Reverse Proxy")
)

$startRow = 28
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Match the post-edit view: scrolled down a bit, selection parked on the
# last entered cell, tab-bar/scrollbar split ratio widened.
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("A35").Select()
$excel.ActiveWindow.TabRatio = 307
